# Add a "Save" column (H) to the sheet, based on whether D (K, strikeouts)
# is above a "save" threshold (~20), mirroring the other 0/1 flag columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the other header cells in row 1 (bold, bordered,
# centered). Copy the formatting from the neighboring "sum" header cell so
# the same style is reused, then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    if ($dVal -gt 20) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
